# "Fruta / hortaliza, semanal" weekly update:
# A brand-new reporting date (2023-09-20, serial 45189) is inserted as a new
# group of 3 rows (Primera/Segunda/Tercera) right above the previously most
# recent historical rows, pushing all subsequent rows down by 3 and growing
# the used range from A1:T122 to A1:T125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 11, shifting the existing rows 11..122 down to 14..125.
$ws.Rows("11:13").Insert()

# Columns that are constant across every data row in this sheet.
$mercadoId  = 1
$mercado    = 'Agrícola del Norte S.A. de Arica'
$region     = 'Arica y Parinacota'
$codreg     = 15
$tipo       = 'Fruta'
$productoId = 100101
$producto   = 'Berries'
$categoriaId = 100112025
$categoria  = 'Frutilla'
$variedad   = 'Sin especificar'
$unidad     = '$/bandeja 3 kilos'
$origen     = 'Región de Arica y Parinacota'
$kgUnidad   = 3

# New data for the inserted rows (Fecha serial 45189 = 2023-09-20).
$nuevasFilas = @(
    @{ Row = 11; Calidad = 'Primera'; Volumen = 80;  PMin = 7000; PMax = 8000; PProm = 7625; PKg = 2542 },
    @{ Row = 12; Calidad = 'Segunda'; Volumen = 120; PMin = 5000; PMax = 6000; PProm = 5417; PKg = 1806 },
    @{ Row = 13; Calidad = 'Tercera'; Volumen = 90;  PMin = 3000; PMax = 4000; PProm = 3556; PKg = 1185 }
)

foreach ($fila in $nuevasFilas) {
    $r = $fila.Row
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = 45189
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $fila.Calidad
    $ws.Cells.Item($r, 13).Value = $fila.Volumen
    $ws.Cells.Item($r, 14).Value = $fila.PMin
    $ws.Cells.Item($r, 15).Value = $fila.PMax
    $ws.Cells.Item($r, 16).Value = $fila.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $fila.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
